# Updates cryptos list data (price and volume changes, plus a few coin re-ordering/replacements)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '42.184.59'
$ws.Range('E2').Value = '  -1.16%  '

# Row 3
$ws.Range('D3').Value = '2.273.20'
$ws.Range('E3').Value = '  -1.34%  '

# Row 4
$ws.Range('E4').Value = '  +0.03%  '

# Row 5
$ws.Range('D5').Value = "'299.13"
$ws.Range('E5').Value = '  -1.49%  '

# Row 6
$ws.Range('D6').Value = "'95.23"
$ws.Range('E6').Value = '  -4.35%  '

# Row 7
$ws.Range('E7').Value = '  -2.60%  '

# Row 8
$ws.Range('E8').Value = '  +0.01%  '

# Row 9
$ws.Range('D9').Value = "'0.490"
$ws.Range('E9').Value = '  -2.84%  '

# Row 10
$ws.Range('D10').Value = "'33.10"
$ws.Range('E10').Value = '  -4.66%  '

# Row 11
$ws.Range('E11').Value = '  -0.60%  '

# Row 12
$ws.Range('D12').Value = "'48.40"
$ws.Range('E12').Value = '  -6.73%  '

# Row 13
$ws.Range('E13').Value = '  +0.82%  '

# Row 14
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = "'15.88"
$ws.Range('E14').Value = '  +1.01%  '

# Row 15
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = "'6.66"
$ws.Range('E15').Value = '  -1.56%  '

# Row 16
$ws.Range('D16').Value = '2.626.00'
$ws.Range('E16').Value = '  -1.41%  '

# Row 17
$ws.Range('D17').Value = '2.262.90'
$ws.Range('E17').Value = '  -1.98%  '

# Row 18
$ws.Range('E18').Value = '  -2.56%  '

# Row 19
$ws.Range('D19').Value = '42.152.69'
$ws.Range('E19').Value = '  -1.06%  '

# Row 20
$ws.Range('D20').Value = "'11.71"
$ws.Range('E20').Value = '  +1.61%  '

# Row 21
$ws.Range('D21').Value = '0.0₃0890'
$ws.Range('E21').Value = '  -1.72%  '

# Row 22
$ws.Range('E22').Value = '  -1.48%  '

# Row 23
$ws.Range('D23').Value = "'66.19"
$ws.Range('E23').Value = '  -2.78%  '

# Row 24
$ws.Range('D24').Value = "'235.15"
$ws.Range('E24').Value = '  -0.17%  '

# Row 25
$ws.Range('E25').Value = '  -1.02%  '

# Row 26
$ws.Range('E26').Value = '  +0.00%  '

# Row 27
$ws.Range('E27').Value = '  -2.88%  '

# Row 28
$ws.Range('D28').Value = "'23.76"
$ws.Range('E28').Value = '  -5.21%  '

# Row 29
$ws.Range('D29').Value = "'2.25"
$ws.Range('E29').Value = '  +2.65%  '

# Row 30
$ws.Range('D30').Value = "'168.31"
$ws.Range('E30').Value = '  +2.91%  '

# Row 31
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = "'33.76"
$ws.Range('E31').Value = '  -2.87%  '

# Row 32
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').Value = "'9.16"
$ws.Range('E32').Value = '  -0.29%  '

# Row 33
$ws.Range('E33').Value = '  +0.03%  '

# Row 34
$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').Value = "'4.64"
$ws.Range('E34').Value = '  +1.43%  '

# Row 35
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').Value = "'4.89"
$ws.Range('E35').Value = '  -2.60%  '

# Row 36
$ws.Range('D36').Value = "'16.76"
$ws.Range('E36').Value = '  -1.12%  '

# Row 37
$ws.Range('D37').Value = "'2.36"
$ws.Range('E37').Value = '  -2.14%  '

# Row 38
$ws.Range('D38').Value = "'0.0688"
$ws.Range('E38').Value = '  -2.77%  '

# Row 39
$ws.Range('E39').Value = '  -2.88%  '

# Row 40
$ws.Range('E40').Value = '  -1.90%  '

# Row 41
$ws.Range('E41').Value = '  -4.71%  '

# Row 42
$ws.Range('E42').Value = '  -2.84%  '

# Row 43
$ws.Range('D43').Value = "'2.29"
$ws.Range('E43').Value = '  -9.07%  '

# Row 44
$ws.Range('D44').Value = '1.961.34'
$ws.Range('E44').Value = '  -0.47%  '

# Row 45
$ws.Range('E45').Value = '  -1.40%  '

# Row 46
$ws.Range('D46').Value = "'17.67"
$ws.Range('E46').Value = '  -4.71%  '

# Row 47
$ws.Range('D47').Value = "'9.60"
$ws.Range('E47').Value = '  -6.32%  '

# Row 48
$ws.Range('E48').Value = '  -4.30%  '

# Row 49
$ws.Range('D49').Value = '2.497.09'
$ws.Range('E49').Value = '  -1.22%  '

# Row 50
$ws.Range('D50').Value = "'52.19"
$ws.Range('E50').Value = '  -5.68%  '

# Row 51
$ws.Range('B51').Value = 'HuobiToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D51').Value = "'2.75"
$ws.Range('E51').Value = '  -2.50%  '
